# Regenerate the handback-status report with the new handoff/handback
# UUID file names, xliff hash, and timestamps (commit: "Generate Report
# for Handback").

$wb = $excel.ActiveWorkbook

$oldUuid1 = "2cc1c537-da89-49fa-a04e-5b5922f03d28"
$oldUuid2 = "b0c6b7d8-42cb-452d-ba11-12a26c399c03"
$newUuid1 = "4a76797a-62b4-4363-a72e-772f7aa5d7e7"
$newUuid2 = "ffff80eb7111-9def-4022-94b6-787d58781a3d"

$newHash  = "4b7ae30990be9f52f1a823ef7d37715c5d818e38"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value = "2016-08-24 00:59:45"

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("G3").Value = "2016-08-24 00:59:45"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\$oldUuid1.md") {
        $hl.TextToDisplay = "e2e\$newUuid1.md"
    } elseif ($hl.TextToDisplay -eq "e2e\$oldUuid2.md") {
        $hl.TextToDisplay = "e2e\$newUuid2.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 00:59:40"
$wsZhCn.Range("K2").Value = "2016-08-24 00:59:56"

$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-24 00:59:40"
$wsZhCn.Range("K3").Value = "2016-08-24 00:59:56"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldUuid1.md") {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldUuid2.md") {
        $hl.TextToDisplay = "$newUuid2.md"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 00:59:45"
$wsDeDe.Range("K2").Value = "2016-08-24 01:00:17"

$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-24 00:59:45"
$wsDeDe.Range("K3").Value = "2016-08-24 01:00:17"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldUuid1.md") {
        $hl.TextToDisplay = "$newUuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldUuid2.md") {
        $hl.TextToDisplay = "$newUuid2.md"
    }
}
